$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix rendering issue: the "LSWR Ch ..." required-reading cells were missing
# the extra trailing space before the literal "<br>" markdown-style line
# break, which broke rendering downstream. Re-set each cell's text with the
# corrected double space before "<br>".
$ws.Range("D2").Value = "LSWR Ch 2 and 3  <br> Light, Singer & Willet 1990, Ch. 2"
$ws.Range("D3").Value = "LSWR Ch 6  <br> [Clayton 2020](https://nautil.us/issue/92/frontiers/how-eugenics-shaped-statistics)"
$ws.Range("D4").Value = "LSWR Ch 11 and 12  <br> [Evans 2020](https://www.newstatesman.com/uncategorized/2020/07/ra-fisher-and-science-hatred)"

# Move/restore the sheet's active selection to D4.
$ws.Range("D4").Select() | Out-Null
